$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.422.61"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "'1.640.08"
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "'303.97"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").Value = "'0.3785"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("D8").Value = "'52.29"
$ws.Range("E8").Value = "  -1.34%  "
$ws.Range("D9").Value = "'0.3642"
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("D10").Value = "'1.248"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("D11").Value = "'0.08103"
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "'22.88"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Value = "'6.632"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "'0.00001252"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "'7.281"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").Value = "'1.639.24"
$ws.Range("E17").Value = "  +2.31%  "
$ws.Range("D18").Value = "'94.10"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "'0.06933"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").Value = "'18.14"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "'6.544"
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("D22").Value = "'1.003"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "'23.430.36"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("D24").Value = "'12.85"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").Value = "'3.255"
$ws.Range("E25").Value = "  +6.01%  "
$ws.Range("D26").Value = "'2.456"
$ws.Range("E26").Value = "  +1.60%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "'149.84"
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("D29").Value = "'5.304"
$ws.Range("E29").Value = "  +0.80%  "
$ws.Range("D30").Value = "'135.99"
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").Value = "'2.323"
$ws.Range("E31").Value = "  -3.50%  "
$ws.Range("D32").Value = "'1.816.99"
$ws.Range("E32").Value = "  +2.08%  "
$ws.Range("D33").Value = "'6.891"
$ws.Range("E33").Value = "  +2.36%  "
$ws.Range("D34").Value = "'10.95"
$ws.Range("E34").Value = "  +6.59%  "
$ws.Range("D35").Value = "'0.9648"
$ws.Range("E35").Value = "  +1.79%  "
$ws.Range("D36").Value = "'0.02866"
$ws.Range("E36").Value = "  +3.76%  "
$ws.Range("D37").Value = "'6.269"
$ws.Range("E37").Value = "  +2.65%  "
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("D39").Value = "'0.07293"
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("D40").Value = "'0.08890"
$ws.Range("E40").Value = "  +1.64%  "
$ws.Range("D41").Value = "'1.373"
$ws.Range("E41").Value = "  -1.85%  "
$ws.Range("D42").Value = "'0.7117"
$ws.Range("E42").Value = "  +0.41%  "
$ws.Range("D43").Value = "'16.44"
$ws.Range("E43").Value = "  +4.31%  "
$ws.Range("D44").Value = "'12.57"
$ws.Range("E44").Value = "  +1.28%  "
$ws.Range("D45").Value = "'0.6562"
$ws.Range("E45").Value = "  +0.74%  "
$ws.Range("D46").Value = "'2.357"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("D47").Value = "'1.002"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").Value = "'3.998"
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("D49").Value = "'0.07995"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("D50").Value = "'1.221"
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("D51").Value = "'127.75"
$ws.Range("E51").Value = "  -4.54%  "
